$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.241.00'
$ws.Cells.Item(2, 5).Value = '  +0.30%  '

$ws.Cells.Item(3, 4).Value = '1.862.36'
$ws.Cells.Item(3, 5).Value = '  -0.05%  '

$ws.Cells.Item(4, 4).Value = '''1.002'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.13%  '

$ws.Cells.Item(5, 4).Value = '''236.89'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.37%  '

$ws.Cells.Item(6, 5).Value = '  +0.11%  '

$ws.Cells.Item(7, 4).Value = '''0.4680'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.39%  '

$ws.Cells.Item(8, 4).Value = '''0.2862'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +1.06%  '

$ws.Cells.Item(9, 4).Value = '''0.06540'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.12%  '

$ws.Cells.Item(10, 4).Value = '''22.02'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +9.25%  '

$ws.Cells.Item(11, 4).Value = '''0.07909'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.17%  '

$ws.Cells.Item(12, 4).Value = '''97.83'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.76%  '

$ws.Cells.Item(13, 4).Value = '1.867.67'
$ws.Cells.Item(13, 5).Value = '  +0.30%  '

$ws.Cells.Item(14, 4).Value = '''5.183'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.22%  '

$ws.Cells.Item(15, 4).Value = '''0.6831'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.12%  '

$ws.Cells.Item(16, 4).Value = '''277.46'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.94%  '

$ws.Cells.Item(17, 4).Value = '30.254.30'
$ws.Cells.Item(17, 5).Value = '  +0.28%  '

$ws.Cells.Item(18, 4).Value = '''13.59'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +7.72%  '

$ws.Cells.Item(19, 5).Value = '  +0.03%  '

$ws.Cells.Item(20, 4).Value = '''0.000007357'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.72%  '

$ws.Cells.Item(21, 4).Value = '2.114.87'
$ws.Cells.Item(21, 5).Value = '  +0.69%  '

$ws.Cells.Item(22, 4).Value = '''5.334'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -2.16%  '

$ws.Cells.Item(23, 4).Value = '''1.002'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.18%  '

$ws.Cells.Item(24, 4).Value = '''6.189'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.79%  '

$ws.Cells.Item(25, 4).Value = '''168.01'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.56%  '

$ws.Cells.Item(26, 4).Value = '''9.242'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.70%  '

$ws.Cells.Item(27, 4).Value = '''19.04'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.79%  '

$ws.Cells.Item(28, 4).Value = '''1.953'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.66%  '

$ws.Cells.Item(29, 4).Value = '''1.388'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +3.01%  '

$ws.Cells.Item(30, 4).Value = '''0.09836'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.29%  '

$ws.Cells.Item(31, 4).Value = '''4.371'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.68%  '

$ws.Cells.Item(32, 4).Value = '''1.486'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.10%  '

$ws.Cells.Item(33, 4).Value = '''4.068'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.11%  '

$ws.Cells.Item(34, 4).Value = '''0.04731'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.54%  '

$ws.Cells.Item(35, 4).Value = '''1.136'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.66%  '

$ws.Cells.Item(36, 4).Value = '''0.7038'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.43%  '

$ws.Cells.Item(37, 4).Value = '''2.711'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.02%  '

$ws.Cells.Item(38, 4).Value = '''0.01877'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.37%  '

$ws.Cells.Item(39, 4).Value = '''2.631'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +4.33%  '

$ws.Cells.Item(40, 4).Value = '''6.278'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.02%  '

$ws.Cells.Item(41, 4).Value = '''75.45'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +4.14%  '

$ws.Cells.Item(42, 4).Value = '''1.953'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.53%  '

$ws.Cells.Item(43, 4).Value = '''0.8517'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.14%  '

$ws.Cells.Item(44, 4).Value = '''0.4170'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.38%  '

$ws.Cells.Item(45, 4).Value = '''1.001'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.06%  '

$ws.Cells.Item(46, 4).Value = '''103.36'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.10%  '

$ws.Cells.Item(47, 4).Value = '''7.203'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.66%  '

$ws.Cells.Item(48, 4).Value = '''956.07'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -3.12%  '

$ws.Cells.Item(49, 4).Value = '''9.236'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.92%  '

$ws.Cells.Item(50, 4).Value = '''34.19'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.13%  '

$ws.Cells.Item(51, 4).Value = '''0.05648'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.11%  '
